$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add three new columns ---
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Copy the formatting (bold/border/center style) from the last existing
# header cell (L1) onto the three newly added header cells.
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- Data rows (rows 2-25): populate the three new columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 13).Value = "after"      # column M
    $ws.Cells.Item($r, 14).Value = 20160636     # column N
    $ws.Cells.Item($r, 15).Value = 3            # column O
}
